$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------------
# 1. Row 6: rename the element from "Extension.valueReference" to the
#    generic "Extension.value[x]" (it is being turned into a sliced
#    value[x] with a dedicated valueReference slice added as row 7).
# -----------------------------------------------------------------------
$ws.Range("B6").Value = ""
$ws.Range("J6").Value = "Reference`n"
$ws.Range("K6").Value = "Value of extension"
$ws.Range("L6").Value = "Value of extension - may be a resource or one of a constrained set of the data types (see Extensibility in the spec for list)."
$ws.Range("AA6").Value = "type:`$this}`n"
$ws.Range("AB6").Value = ""
$ws.Range("AD6").Value = "closed"
$ws.Range("AE6").Value = "Extension.value[x]"

# -----------------------------------------------------------------------
# 2. Add new row 7 - the "valueReference" slice of Extension.value[x],
#    carrying the content that used to live on row 6.
#    Copy row 6's formatting first so every cell in the new row keeps
#    the same style, then fill in the text that differs from blank.
# -----------------------------------------------------------------------
$ws.Range("A6:AJ6").Copy()
$ws.Range("A7:AJ7").PasteSpecial(-4122)
$ws.Rows.Item(7).Hidden = $true

$ws.Range("A7").Value = "Extension.value[x]"
$ws.Range("B7").Value = "valueReference"
$ws.Range("E7").Value = "0"
$ws.Range("F7").Value = "1"
$ws.Range("J7").Value = "Reference(https://fhir.hl7.org.uk/STU3/StructureDefinition/CareConnect-Organization-1)`n"
$ws.Range("K7").Value = "The patient's nominated pharmacy"
$ws.Range("L7").Value = "The patient's nominated pharmacy."
$ws.Range("AE7").Value = "Extension.value[x]"
$ws.Range("AF7").Value = "0"
$ws.Range("AG7").Value = "1"
$ws.Range("AJ7").Value = "N/A"

# -----------------------------------------------------------------------
# 3. Column A width shrinks now that "Extension.value[x]" (the bestFit
#    label) is shorter than "Extension.valueReference".
# -----------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 18.15

# -----------------------------------------------------------------------
# 4. Extend the AutoFilter range down to the new row and re-apply the
#    same filter criteria that were already active (toggling the filter
#    off then on again is required to change its range via COM).
# -----------------------------------------------------------------------
$ws.Range("A1:AJ6").AutoFilter()
$ws.Range("A1:AJ7").AutoFilter(7, "<> ", 1)
$ws.Range("A1:AJ7").AutoFilter(27, @(""), 7)

# -----------------------------------------------------------------------
# 5. Conditional formatting covered rows 2-5; it now needs to cover the
#    re-labelled row 6 as well (rows 2-6).
# -----------------------------------------------------------------------
$fc = $ws.Range("A2:AI5").FormatConditions.Item(1)
$fc.ModifyAppliesToRange($ws.Range("A2:AI6"))

# -----------------------------------------------------------------------
# 6. The hidden _FilterDatabase defined name must track the new range.
# -----------------------------------------------------------------------
$wb.Names.Item(1).RefersTo = "=Elements!`$A`$1:`$AJ`$7"
